# Update column G ("K") values on Sheet1 to reflect the regenerated
# save_data (using K instead of Strike#, recalculated std/mean, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 3
    3  = 1
    4  = 0
    5  = 0
    6  = 1
    7  = 2
    8  = 2
    9  = 0
    10 = 3
    11 = 0
    12 = 1
    13 = 1
    14 = 2
    15 = 2
    16 = 0
    17 = 4
    18 = 1
    19 = 2
    20 = 1
    21 = 1
    22 = 0
    23 = 1
    24 = 4
    25 = 3
    26 = 5
    27 = 3
    28 = 4
    29 = 0
    30 = 2
    31 = 0
    32 = 3
    34 = 3
    36 = 3
    37 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
